$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add()
$ws.Name = "Sheet1"
$old = $wb.Worksheets.Item("DOCA Installation")
$old.Delete()

# Header row
$ws.Range("A1").Value = '#'
$ws.Range("B1").Value = 'Command'
$ws.Range("C1").Value = 'Output'

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 'sudo –i'
$ws.Range("C2").Value = '[sudo] password for student:
Sorry, try again.
[sudo] password for student:
sudo: –i​: command not found'

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 'cd /home/student/AI_Infra/module5/hands_on_1'
$ws.Range("C3").Value = ''

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = '/usr/sbin/ofed_uninstall.sh --force'
$ws.Range("C4").Value = 'bash: /usr/sbin/ofed_uninstall.sh: No such file or directory'

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 'sudo apt-get autoremove'
$ws.Range("C5").Value = 'Reading package lists... Done
Building dependency tree... Done
Reading state information... Done
0 upgraded, 0 newly installed, 0 to remove and 327 not upgraded.'

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 'dpkg -i doca-host_3.1.0-091000-25.07-ubuntu2204_amd64.deb'
$ws.Range("C6").Value = 'dpkg: error: requested operation requires superuser privilege'

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 'sudo -i'
$ws.Range("C7").Value = 'root@acad14:~#'

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 'for f in $( dpkg --list | grep doca | awk ''{print $2}'' ); do echo $f ; apt remove --purge $f -y ; done'
$ws.Range("C8").Value = 'doca-host
(Reading database ... 194224 files and directories currently installed.)
Removing doca-host (3.1.0-091000-25.07-ubuntu2204) ...
Purging configuration files for doca-host (3.1.0-091000-25.07-ubuntu2204) ...'

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = '/usr/sbin/ofed_uninstall.sh --force'
$ws.Range("C9").Value = '-bash: /usr/sbin/ofed_uninstall.sh: No such file or directory'

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 'cd /home/student/AI_Infra/module5/hands_on_1'
$ws.Range("C10").Value = ''

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 'dpkg -i doca-host_3.1.0-091000-25.07-ubuntu2204_amd64.deb'
$ws.Range("C11").Value = 'Selecting previously unselected package doca-host.
(Reading database ... 194224 files and directories currently installed.)
Preparing to unpack doca-host_3.1.0-091000-25.07-ubuntu2204_amd64.deb ...
Unpacking doca-host (3.1.0-091000-25.07-ubuntu2204) ...
Setting up doca-host (3.1.0-091000-25.07-ubuntu2204) ...'

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 'apt update'
$ws.Range("C12").Value = 'Fetched 21.9 MB in 4s (6,219 kB/s)
Reading package lists... Done
Building dependency tree... Done
Reading state information... Done
347 packages can be upgraded. Run ''apt list --upgradable'' to see them.'

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 'systemctl status rshim'
$ws.Range("C13").Value = 'Unit rshim.service could not be found.'

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 'apt install -y rshim'
$ws.Range("C14").Value = 'Reading package lists... Done
Building dependency tree... Done
Reading state information... Done
The following NEW packages will be installed:
  rshim
0 upgraded, 1 newly installed, 0 to remove and 347 not upgraded.
Need to get 14.2 kB of archives.
After this operation, 45.1 kB of additional disk space will be used.
Get:1 http://archive.ubuntu.com/ubuntu jammy/main amd64 rshim amd64 2.0.6-0ubuntu1 [14.2 kB]
Fetched 14.2 kB in 0s (35.6 kB/s)
Selecting previously unselected package rshim.
(Reading database ... 194224 files and directories currently installed.)
Preparing to unpack .../rshim_2.0.6-0ubuntu1_amd64.deb ...
Unpacking rshim (2.0.6-0ubuntu1) ...
Setting up rshim (2.0.6-0ubuntu1) ...
Created symlink /etc/systemd/system/multi-user.target.wants/rshim.service → /lib/systemd/system/rshim.service.'

$ws.Range("A15").Value = 14
$ws.Range("B15").Value = 'sudo systemctl start rshim'
$ws.Range("C15").Value = ''

$ws.Range("A16").Value = 15
$ws.Range("B16").Value = 'apt install -y doca-all mlnx-fw-updater'
$ws.Range("C16").Value = '21 upgraded, 160 newly installed, 7 to remove and 325 not upgraded.
Need to get 14.0 MB/416 MB of archives.
After this operation, 1,315 MB of additional disk space will be used.

Device #1:
----------
  Device Type:      ConnectX7
  Part Number:      MCX755106AC-HEA_Ax
  FW:               28.46.1006
  Status:           Up to date

DKMS builds: xpmem, kernel-mft, knem, mlnx-ofed-kernel, srp, iser, isert

Setting up doca-all (3.1.0-091000) ...'

$ws.Range("A17").Value = 16
$ws.Range("B17").Value = 'exit'
$ws.Range("C17").Value = 'logout'

$ws.Range("A18").Value = 17
$ws.Range("B18").Value = 'sudo minicom -D /dev/rshim0/console'
$ws.Range("C18").Value = 'minicom: cannot open /dev/rshim0/console: No such file or directory'

$ws.Range("A19").Value = 18
$ws.Range("B19").Value = 'sudo systemctl start rshim'
$ws.Range("C19").Value = ''

$ws.Range("A20").Value = 19
$ws.Range("B20").Value = 'sudo minicom -D /dev/rshim0/console'
$ws.Range("C20").Value = 'minicom: cannot open /dev/rshim0/console: No such file or directory'

$ws.Range("A21").Value = 20
$ws.Range("B21").Value = 'sudo systemctl status rshim'
$ws.Range("C21").Value = '● rshim.service - rshim driver for BlueField SoC
     Active: active (running) since Sun 2025-11-23 07:28:43 PST
   Main PID: 8968 (rshim)
rshim[8968]: Probing pcie-0000:a0:00.2(uio)
rshim[8968]: Create rshim pcie-0000:a0:00.2
rshim[8968]: another backend already attached
rshim[8968]: rshim0 entering drop mode'

$ws.Range("A22").Value = 21
$ws.Range("B22").Value = 'sudo systemctl stop rshim'
$ws.Range("C22").Value = ''

$ws.Range("A23").Value = 22
$ws.Range("B23").Value = 'sudo systemctl start rshim'
$ws.Range("C23").Value = ''

$ws.Range("A24").Value = 23
$ws.Range("B24").Value = 'sudo minicom -D /dev/rshim0/console'
$ws.Range("C24").Value = 'Welcome to minicom 2.8

Port /dev/rshim0/console, 07:33:52

Press CTRL-A Z for help on special keys

acad14-DPU login:'

$ws.Range("A25").Value = 24
$ws.Range("B25").Value = 'ubuntu'
$ws.Range("C25").Value = 'Password:'

$ws.Range("A26").Value = 27
$ws.Range("B26").Value = 'sudo -i'
$ws.Range("C26").Value = 'sudo: unable to resolve host acad14-DPU: Temporary failure in name resolution
root@acad14-DPU:~#'

$ws.Range("A27").Value = 28
$ws.Range("B27").Value = 'sudo bfver'
$ws.Range("C27").Value = '--/dev/mmcblk0boot0
BlueField ATF version: v2.2(release):4.9.0-25-g0ce57e322
BlueField UEFI version: 4.9.0-46-g7e3911bd4d
BlueField BSP version: 4.9.0.13378

OS Release Version: bf-bundle-2.9.0-90_24.10_ubuntu-22.04_prod'

# Header styling: bold font, borders, centered/top aligned (matches the
# "apt install -y rshim" era template: only the header row carries an
# explicit style now, data rows are left at the default style).
$headerRange = $ws.Range("A1:C1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Multi-line cell text triggers simulated auto row-height; restore the
# natural (non-custom) row heights like the original plain template had.
$ws.Rows("1:27").AutoFit()

$ws.Range("A1").Select()
